$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 64
$ws.Range("H64").Value = 7569.2104
$ws.Range("I64").Value = 5105.8335
$ws.Range("J64").Value = 8706.154
$ws.Range("K64").Value = 5105.8335
$ws.Range("L64").Value = 8706.154
$ws.Range("M64").Value = -4857.8335
$ws.Range("N64").Value = -9202.154

# Row 67
$ws.Range("H67").Value = 7569.2104
$ws.Range("I67").Value = 5105.8335
$ws.Range("J67").Value = 8706.154
$ws.Range("K67").Value = 5105.8335
$ws.Range("L67").Value = 8706.154
$ws.Range("M67").Value = -4247.8335
$ws.Range("N67").Value = -10422.154

# Row 98
$ws.Range("H98").Value = 1589.7037
$ws.Range("I98").Value = 1616.7084
$ws.Range("J98").Value = 1373.6666
$ws.Range("K98").Value = 1616.7084
$ws.Range("L98").Value = 1373.6666
$ws.Range("M98").Value = -118.7084
$ws.Range("N98").Value = -4369.6666

# Row 115
$ws.Range("H115").Value = 638.75
$ws.Range("I115").Value = 385
$ws.Range("J115").Value = 1400
$ws.Range("K115").Value = 1155
$ws.Range("L115").Value = 4200
$ws.Range("M115").Value = 412
$ws.Range("N115").Value = -7334

# Row 122
$ws.Range("H122").Value = 1589.7037
$ws.Range("I122").Value = 1616.7084
$ws.Range("J122").Value = 1373.6666
$ws.Range("K122").Value = 4850.1252
$ws.Range("L122").Value = 4120.9998
$ws.Range("M122").Value = -2400.1252
$ws.Range("N122").Value = -9020.9998

# Row 129
$ws.Range("H129").Value = 1780
$ws.Range("I129").Value = 1438.4615
$ws.Range("K129").Value = 4315.3845
$ws.Range("M129").Value = 684.6154999999999

# Row 132
$ws.Range("H132").Value = 71058.82000000001
$ws.Range("I132").Value = 74718.59
$ws.Range("K132").Value = 224155.77
$ws.Range("M132").Value = -221625.77

# Row 135
$ws.Range("H135").Value = 1029.5143
$ws.Range("I135").Value = 813.9545000000001
$ws.Range("J135").Value = 1394.3077
$ws.Range("K135").Value = 7325.5905
$ws.Range("L135").Value = 12548.7693
$ws.Range("M135").Value = -4790.5905
$ws.Range("N135").Value = -17618.7693

# Row 137
$ws.Range("H137").Value = 1807188.2
$ws.Range("I137").Value = 3738.4443
$ws.Range("J137").Value = 3610638
$ws.Range("K137").Value = 11215.3329
$ws.Range("L137").Value = 10831914
$ws.Range("M137").Value = -8665.332900000001
$ws.Range("N137").Value = -10837014

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 8918.263000000001
$ws.Range("I2").Value = 1609.7858
$ws.Range("J2").Value = 29382
$ws.Range("K2").Value = 1609.7858
$ws.Range("L2").Value = 29382
$ws.Range("M2").Value = -1496.7858
$ws.Range("N2").Value = -29608

# Row 32
$ws.Range("H32").Value = 21703.963
$ws.Range("I32").Value = 22113.846
$ws.Range("K32").Value = 22113.846
$ws.Range("M32").Value = -21826.846

# Row 61
$ws.Range("H61").Value = 3704547.5
$ws.Range("I61").Value = 3704547.5
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 3704547.5
$ws.Range("L61").Value = 0
$ws.Range("M61").ClearContents()
$ws.Range("N61").Value = -3704335.5

# Row 74
$ws.Range("H74").Value = 3540.8845
$ws.Range("I74").Value = 1072.8
$ws.Range("K74").Value = 1072.8
$ws.Range("M74").Value = -198.8

# Row 77
$ws.Range("H77").Value = 3540.8845
$ws.Range("I77").Value = 1072.8
$ws.Range("K77").Value = 5364
$ws.Range("M77").Value = -996

# Row 116
$ws.Range("H116").Value = 8918.263000000001
$ws.Range("I116").Value = 1609.7858
$ws.Range("J116").Value = 29382
$ws.Range("K116").Value = 1609.7858
$ws.Range("L116").Value = 29382
$ws.Range("M116").Value = 684.2141999999999
$ws.Range("N116").Value = -33970

# Row 132
$ws.Range("H132").Value = 410331.6
$ws.Range("I132").Value = 466889.8
$ws.Range("J132").Value = 4997.6665
$ws.Range("K132").Value = 1400669.4
$ws.Range("L132").Value = 14992.9995
$ws.Range("M132").Value = -1398139.4
$ws.Range("N132").Value = -20052.9995

# Row 136
$ws.Range("H136").Value = 3704547.5
$ws.Range("I136").Value = 3704547.5
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 11113642.5
$ws.Range("L136").Value = 0
$ws.Range("M136").ClearContents()
$ws.Range("N136").Value = -11111092.5

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 8918.263000000001
$ws.Range("I3").Value = 1609.7858
$ws.Range("J3").Value = 29382
$ws.Range("K3").Value = 1609.7858
$ws.Range("L3").Value = 29382
$ws.Range("M3").Value = -1495.7858
$ws.Range("N3").Value = -29610

# Row 86
$ws.Range("H86").Value = 1210.4286
$ws.Range("I86").Value = 1223.1818
$ws.Range("K86").Value = 1223.1818
$ws.Range("M86").Value = -100.1818000000001

# Row 89
$ws.Range("H89").Value = 1210.4286
$ws.Range("I89").Value = 1223.1818
$ws.Range("K89").Value = 6115.909000000001
$ws.Range("M89").Value = -499.9090000000006

$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Range("H22").Value = 713.86664
$ws.Range("I22").Value = 726.2308
$ws.Range("J22").Value = 633.5
$ws.Range("K22").Value = 726.2308
$ws.Range("L22").Value = 633.5
$ws.Range("M22").Value = -376.2308
$ws.Range("N22").Value = -1333.5

# Row 31
$ws.Range("H31").Value = 20300.166
$ws.Range("I31").Value = 10418.2
$ws.Range("K31").Value = 10418.2
$ws.Range("M31").Value = -10123.2

# Row 34
$ws.Range("H34").Value = 20300.166
$ws.Range("I34").Value = 10418.2
$ws.Range("K34").Value = 10418.2
$ws.Range("M34").Value = -10216.2

# Row 51
$ws.Range("H51").Value = 15848.333
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").ClearContents()

# Row 59
$ws.Range("H59").Value = 62249
$ws.Range("I59").Value = 4500
$ws.Range("J59").Value = 119998
$ws.Range("K59").Value = 4500
$ws.Range("L59").Value = 119998
$ws.Range("M59").Value = -3355
$ws.Range("N59").Value = -122288

# Row 60
$ws.Range("H60").Value = 21124.75
$ws.Range("I60").Value = 21124.75
$ws.Range("J60").Value = 0
$ws.Range("K60").Value = 21124.75
$ws.Range("L60").Value = 0
$ws.Range("M60").ClearContents()
$ws.Range("N60").Value = -20613.75

# Row 61
$ws.Range("H61").Value = 15848.333
$ws.Range("J61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("N61").ClearContents()

# Row 134
$ws.Range("H134").Value = 5475.943
$ws.Range("I134").Value = 5771.6875
$ws.Range("J134").Value = 2321.3333
$ws.Range("K134").Value = 17315.0625
$ws.Range("L134").Value = 6963.999899999999
$ws.Range("M134").Value = -14780.0625
$ws.Range("N134").Value = -12033.9999

$ws = $wb.Worksheets.Item("CUL")
# Row 4
$ws.Range("H4").Value = 115677816
$ws.Range("I4").Value = 115677816
$ws.Range("K4").Value = 347033448
$ws.Range("M4").Value = -347033336

$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 4731.8335
$ws.Range("I70").Value = 4747.75
$ws.Range("J70").Value = 4700
$ws.Range("K70").Value = 4747.75
$ws.Range("L70").Value = 4700
$ws.Range("M70").Value = -4477.75
$ws.Range("N70").Value = -5240

# Row 73
$ws.Range("H73").Value = 4731.8335
$ws.Range("I73").Value = 4747.75
$ws.Range("J73").Value = 4700
$ws.Range("K73").Value = 4747.75
$ws.Range("L73").Value = 4700
$ws.Range("M73").Value = -3811.75
$ws.Range("N73").Value = -6572

# Row 92
$ws.Range("H92").Value = 68100.39999999999
$ws.Range("J92").Value = 68100.39999999999
$ws.Range("L92").Value = 68100.39999999999
$ws.Range("N92").Value = -71844.39999999999

# Row 98
$ws.Range("H98").Value = 75625.39999999999
$ws.Range("J98").Value = 75625.39999999999
$ws.Range("L98").Value = 75625.39999999999
$ws.Range("N98").Value = -81615.39999999999

# Row 102
$ws.Range("H102").Value = 2818.5264
$ws.Range("J102").Value = 3834.6428
$ws.Range("L102").Value = 3834.6428
$ws.Range("N102").Value = -7078.6428

# Row 123
$ws.Range("H123").Value = 69853
$ws.Range("J123").Value = 69853
$ws.Range("L123").Value = 69853
$ws.Range("N123").Value = -74753

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 3651.16
$ws.Range("I7").Value = 3363.75
$ws.Range("J7").Value = 4800.8
$ws.Range("K7").Value = 3363.75
$ws.Range("L7").Value = 4800.8
$ws.Range("M7").Value = -3251.75
$ws.Range("N7").Value = -5024.8

# Row 61
$ws.Range("H61").Value = 2528.2
$ws.Range("I61").Value = 1624.5
$ws.Range("K61").Value = 1624.5
$ws.Range("M61").Value = -1422.5

# Row 113
$ws.Range("H113").Value = 2528.2
$ws.Range("I113").Value = 1624.5
$ws.Range("K113").Value = 1624.5
$ws.Range("M113").Value = 545.5

# Row 126
$ws.Range("H126").Value = 3651.16
$ws.Range("I126").Value = 3363.75
$ws.Range("J126").Value = 4800.8
$ws.Range("K126").Value = 10091.25
$ws.Range("L126").Value = 14402.4
$ws.Range("M126").Value = -7621.25
$ws.Range("N126").Value = -19342.4

# Row 132
$ws.Range("H132").Value = 3486868.5
$ws.Range("I132").Value = 4356085.5
$ws.Range("K132").Value = 13068256.5
$ws.Range("M132").Value = -13065726.5

$ws = $wb.Worksheets.Item("WVR")
# Row 103
$ws.Range("H103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("L103").ClearContents()
$ws.Range("N103").Value = 0

# Row 132
$ws.Range("H132").Value = 6711119
$ws.Range("I132").Value = 7456047
$ws.Range("K132").Value = 22368141
$ws.Range("M132").Value = -22365611

# Row 136
$ws.Range("H136").Value = 10006.928
$ws.Range("I136").Value = 11941.113
$ws.Range("K136").Value = 35823.339
$ws.Range("M136").Value = -33273.339
